$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-25 04:18:35"
$ws.Range("N2").Value = "1.0 °C 3:50 TU"
$ws.Range("E3").Value = "2026-02-25 04:18:38"
$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = "31%"
$ws.Range("E4").Value = "2026-02-25 04:18:40"
$ws.Range("N4").Value = "1.8 °C 3:44 TU"
$ws.Range("O4").Value = "2.9 °C"
$ws.Range("E5").Value = "2026-02-25 04:18:43"
$ws.Range("K5").Value = "-0.1 MJ/m2"
$ws.Range("O5").Value = "5.3 °C"
$ws.Range("E6").Value = "2026-02-25 04:18:45"
$ws.Range("J6").Value = "1018.9 hPa"
$ws.Range("N6").Value = "7.6 °C 3:50 TU"
$ws.Range("O6").Value = "8.9 °C"
$ws.Range("E7").Value = "2026-02-25 04:18:48"
$ws.Range("L7").Value = "9.4 km/h - 86º 3:52 TU"
$ws.Range("O7").Value = "11.1 °C"
$ws.Range("E8").Value = "2026-02-25 04:18:50"
$ws.Range("N8").Value = "14.9 °C 3:34 TU"
$ws.Range("E9").Value = "2026-02-25 04:18:53"
$ws.Range("N9").Value = "3.4 °C 3:37 TU"
$ws.Range("O9").Value = "5.5 °C"
$ws.Range("E10").Value = "2026-02-25 04:18:55"
$ws.Range("N10").Value = "2.5 °C 3:55 TU"
$ws.Range("O10").Value = "4.1 °C"
$ws.Range("E11").Value = "2026-02-25 04:18:58"
$ws.Range("N11").Value = "2.1 °C 3:45 TU"
$ws.Range("O11").Value = "3.2 °C"
$ws.Range("E12").Value = "2026-02-25 04:19:00"
$ws.Range("O12").Value = "5.2 °C"
$ws.Range("E13").Value = "2026-02-25 04:19:03"
$ws.Range("J13").Value = "1026.5 hPa"
$ws.Range("N13").Value = "-2.8 °C 3:39 TU"
$ws.Range("O13").Value = "-1.3 °C"
$ws.Range("E14").Value = "2026-02-25 04:19:05"
$ws.Range("N14").Value = "4.2 °C 3:47 TU"
$ws.Range("O14").Value = "5.8 °C"
$ws.Range("E15").Value = "2026-02-25 04:19:08"
$ws.Range("E16").Value = "2026-02-25 04:19:10"
$ws.Range("H16").NumberFormat = "@"
$ws.Range("H16").Value = "16%"
$ws.Range("O16").Value = "3.7 °C"
$ws.Range("E17").Value = "2026-02-25 04:19:12"
$ws.Range("N17").Value = "7.9 °C 3:30 TU"
$ws.Range("O17").Value = "9.3 °C"
$ws.Range("E18").Value = "2026-02-25 04:19:15"
$ws.Range("N18").Value = "5.0 °C 3:59 TU"
$ws.Range("O18").Value = "6.5 °C"
$ws.Range("E19").Value = "2026-02-25 04:19:18"
$ws.Range("H19").NumberFormat = "@"
$ws.Range("H19").Value = "58%"
$ws.Range("E20").Value = "2026-02-25 04:19:20"
$ws.Range("H20").NumberFormat = "@"
$ws.Range("H20").Value = "44%"
$ws.Range("O20").Value = "2.7 °C"
$ws.Range("E21").Value = "2026-02-25 04:19:23"
$ws.Range("H21").NumberFormat = "@"
$ws.Range("H21").Value = "74%"
$ws.Range("J21").Value = "1023.3 hPa"
$ws.Range("N21").Value = "2.3 °C 3:50 TU"
$ws.Range("O21").Value = "3.9 °C"
$ws.Range("E22").Value = "2026-02-25 04:19:25"
$ws.Range("E23").Value = "2026-02-25 04:19:28"
$ws.Range("H23").NumberFormat = "@"
$ws.Range("H23").Value = "28%"
$ws.Range("L23").Value = "10.8 km/h - 0º 3:39 TU"
$ws.Range("E24").Value = "2026-02-25 04:19:31"
$ws.Range("O24").Value = "3.9 °C"
$ws.Range("E25").Value = "2026-02-25 04:19:33"
$ws.Range("N25").Value = "2.2 °C 3:36 TU"
$ws.Range("O25").Value = "3.3 °C"
$ws.Range("E26").Value = "2026-02-25 04:19:35"
$ws.Range("H26").NumberFormat = "@"
$ws.Range("H26").Value = "40%"
$ws.Range("J26").Value = "1018.8 hPa"
$ws.Range("K26").Value = "-0.1 MJ/m2"
$ws.Range("E27").Value = "2026-02-25 04:19:37"
$ws.Range("L27").Value = "24.1 km/h - 243º 3:56 TU"
$ws.Range("M27").Value = "4.8 °C 3:50 TU"
$ws.Range("E28").Value = "2026-02-25 04:19:40"
$ws.Range("J28").Value = "1020.3 hPa"
$ws.Range("N28").Value = "2.7 °C 3:59 TU"
$ws.Range("O28").Value = "4.0 °C"
$ws.Range("E29").Value = "2026-02-25 04:19:43"
$ws.Range("O29").Value = "8.9 °C"
$ws.Range("E30").Value = "2026-02-25 04:19:45"
$ws.Range("J30").Value = "1019.2 hPa"
$ws.Range("N30").Value = "6.9 °C 3:56 TU"
$ws.Range("O30").Value = "7.8 °C"
$ws.Range("E31").Value = "2026-02-25 04:19:48"
$ws.Range("H31").NumberFormat = "@"
$ws.Range("H31").Value = "91%"
$ws.Range("O31").Value = "10.3 °C"
$ws.Range("E32").Value = "2026-02-25 04:19:51"
$ws.Range("H32").NumberFormat = "@"
$ws.Range("H32").Value = "72%"
$ws.Range("E33").Value = "2026-02-25 04:19:53"
$ws.Range("J33").Value = "1023.9 hPa"
$ws.Range("N33").Value = "1.1 °C 3:57 TU"
$ws.Range("O33").Value = "2.5 °C"
$ws.Range("E34").Value = "2026-02-25 04:19:56"
$ws.Range("N34").Value = "0.3 °C 3:42 TU"
$ws.Range("O34").Value = "1.7 °C"
$ws.Range("E35").Value = "2026-02-25 04:19:59"
$ws.Range("H35").NumberFormat = "@"
$ws.Range("H35").Value = "34%"
$ws.Range("O35").Value = "10.2 °C"
$ws.Range("E36").Value = "2026-02-25 04:20:01"
$ws.Range("J36").Value = "1019.0 hPa"
$ws.Range("L36").Value = "20.9 km/h - 19º 3:58 TU"
$ws.Range("M36").Value = "10.0 °C 3:59 TU"
$ws.Range("O36").Value = "8.2 °C"
$ws.Range("E37").Value = "2026-02-25 04:20:04"
$ws.Range("N37").Value = "0.2 °C 3:55 TU"
$ws.Range("O37").Value = "1.1 °C"
$ws.Range("E38").Value = "2026-02-25 04:20:06"
$ws.Range("L38").Value = "10.1 km/h - 313º 3:44 TU"
$ws.Range("E39").Value = "2026-02-25 04:20:09"
$ws.Range("E40").Value = "2026-02-25 04:20:11"
$ws.Range("N40").Value = "0.6 °C 3:45 TU"
$ws.Range("O40").Value = "1.8 °C"
$ws.Range("E41").Value = "2026-02-25 04:20:14"
$ws.Range("I41").Value = "0.1 mm"
$ws.Range("O41").Value = "8.8 °C"
$ws.Range("E42").Value = "2026-02-25 04:20:17"
$ws.Range("E43").Value = "2026-02-25 04:20:19"
$ws.Range("H43").NumberFormat = "@"
$ws.Range("H43").Value = "91%"
$ws.Range("O43").Value = "3.9 °C"
$ws.Range("E44").Value = "2026-02-25 04:20:22"
$ws.Range("H44").NumberFormat = "@"
$ws.Range("H44").Value = "45%"
$ws.Range("O44").Value = "-0.4 °C"
$ws.Range("E45").Value = "2026-02-25 04:20:24"
$ws.Range("E46").Value = "2026-02-25 04:20:27"
$ws.Range("N46").Value = "2.6 °C 3:58 TU"
$ws.Range("O46").Value = "3.8 °C"
